# Reorder the header row (row 1) on the two "BiomedicalConcept" sheets so
# that packageDate/packageType move from columns D/E to the front (A/B),
# pushing the BiomedicalConcept_* columns to C/D/E, and swapping the order
# of categories/parentConceptId (now parentConceptId=F, categories=G).
#
# Old order: BiomedicalConcept_conceptId, BiomedicalConcept_ncitCode,
#            BiomedicalConcept_href, packageDate, packageType, categories,
#            parentConceptId, ...
# New order: packageDate, packageType, BiomedicalConcept_conceptId,
#            BiomedicalConcept_ncitCode, BiomedicalConcept_href,
#            parentConceptId, categories, ...

$wb = $excel.ActiveWorkbook

$newHeaders = @(
    "packageDate",
    "packageType",
    "BiomedicalConcept_conceptId",
    "BiomedicalConcept_ncitCode",
    "BiomedicalConcept_href",
    "parentConceptId",
    "categories"
)

$sheetNames = @("BiomedicalConcept", "BiomedicalConcept1")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($i = 0; $i -lt $newHeaders.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
    }
}

# The "packageType" data validation (list "bc") followed the packageType
# column, which moved from E to B.
$wsMain = $wb.Worksheets.Item("BiomedicalConcept")
$wsMain.Range("E2:E1048576").Validation.Delete()
$wsMain.Range("B2:B1048576").Validation.Add(3, 1, 1, '"bc"')
